$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold numeric-looking text (e.g. "30.535.46", "1.001").
# Force a Text number format first so the COM layer stores the literal string
# instead of silently coercing it to a Double, then drop back to the default
# "Normal" style afterwards so no stray number-format style lingers on the cells.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.557.81"
$ws.Range("E2").Value = "  -0.47%  "

# Row 3
$ws.Range("D3").Value = "1.874.68"
$ws.Range("E3").Value = "  -0.90%  "

# Row 4
$ws.Range("D4").Value = "1.000"

# Row 5
$ws.Range("D5").Value = "247.56"
$ws.Range("E5").Value = "  +0.53%  "

# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").Value = "0.4759"
$ws.Range("E7").Value = "  -0.68%  "

# Row 8
$ws.Range("D8").Value = "0.2912"
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").Value = "0.06502"
$ws.Range("E9").Value = "  -1.15%  "

# Row 10
$ws.Range("D10").Value = "21.85"
$ws.Range("E10").Value = "  +0.43%  "

# Row 11
$ws.Range("D11").Value = "0.07755"
$ws.Range("E11").Value = "  -0.42%  "

# Row 12
$ws.Range("D12").Value = "0.7382"
$ws.Range("E12").Value = "  -0.61%  "

# Row 13
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "96.13"
$ws.Range("E13").Value = "  -1.51%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.872.38"
$ws.Range("E14").Value = "  -1.12%  "

# Row 15
$ws.Range("D15").Value = "5.174"
$ws.Range("E15").Value = "  -0.46%  "

# Row 16
$ws.Range("D16").Value = "274.38"
$ws.Range("E16").Value = "  -2.57%  "

# Row 17
$ws.Range("D17").Value = "30.598.78"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("E18").Value = "  -2.64%  "

# Row 19
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("D20").Value = "0.000007505"
$ws.Range("E20").Value = "  -1.56%  "

# Row 21
$ws.Range("D21").Value = "2.118.04"
$ws.Range("E21").Value = "  -1.33%  "

# Row 22
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").Value = "5.220"
$ws.Range("E23").Value = "  -1.44%  "

# Row 24
$ws.Range("D24").Value = "6.170"

# Row 25
$ws.Range("D25").Value = "9.175"
$ws.Range("E25").Value = "  -2.11%  "

# Row 26
$ws.Range("D26").Value = "164.84"
$ws.Range("E26").Value = "  -0.67%  "

# Row 27
$ws.Range("D27").Value = "18.81"
$ws.Range("E27").Value = "  -1.80%  "

# Row 28
$ws.Range("D28").Value = "1.909"
$ws.Range("E28").Value = "  -3.14%  "

# Row 29
$ws.Range("D29").Value = "0.09850"
$ws.Range("E29").Value = "  -1.34%  "

# Row 30
$ws.Range("D30").Value = "1.337"
$ws.Range("E30").Value = "  -2.52%  "

# Row 31
$ws.Range("D31").Value = "1.496"
$ws.Range("E31").Value = "  -1.59%  "

# Row 32
$ws.Range("D32").Value = "4.248"
$ws.Range("E32").Value = "  -2.78%  "

# Row 33
$ws.Range("D33").Value = "4.082"
$ws.Range("E33").Value = "  -1.40%  "

# Row 34
$ws.Range("D34").Value = "0.04789"
$ws.Range("E34").Value = "  +0.10%  "

# Row 35
$ws.Range("D35").Value = "1.121"
$ws.Range("E35").Value = "  -0.98%  "

# Row 36
$ws.Range("D36").Value = "0.6935"
$ws.Range("E36").Value = "  -1.77%  "

# Row 37
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +0.03%  "

# Row 38
$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  -1.16%  "

# Row 39
$ws.Range("D39").Value = "2.757"
$ws.Range("E39").Value = "  -0.37%  "

# Row 40
$ws.Range("D40").Value = "6.277"
$ws.Range("E40").Value = "  -2.20%  "

# Row 41
$ws.Range("D41").Value = "73.32"
$ws.Range("E41").Value = "  +3.71%  "

# Row 42
$ws.Range("D42").Value = "1.982"
$ws.Range("E42").Value = "  +2.44%  "

# Row 43
$ws.Range("D43").Value = "0.4202"
$ws.Range("E43").Value = "  -0.64%  "

# Row 44
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("D45").Value = "0.8344"
$ws.Range("E45").Value = "  -1.52%  "

# Row 46
$ws.Range("D46").Value = "101.48"
$ws.Range("E46").Value = "  -1.03%  "

# Row 47
$ws.Range("D47").Value = "9.403"
$ws.Range("E47").Value = "  +0.48%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "6.965"
$ws.Range("E48").Value = "  -2.90%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "35.26"
$ws.Range("E49").Value = "  -0.46%  "

# Row 50
$ws.Range("D50").Value = "912.94"
$ws.Range("E50").Value = "  -2.64%  "

# Row 51
$ws.Range("D51").Value = "0.05667"
$ws.Range("E51").Value = "  +0.97%  "

# Restore the default style now that the literal text values are committed.
$priceRange.Style = "Normal"
